# Update Data Sources from LFX (#91)
# Re-applies the table style used by the "Member Benefits" tables across
# the deck: {48338D13-A9CA-4A3B-89B7-9073ACB4FE37} -> {631E9890-4A1E-43CC-99DB-8398D616D4D7}

$oldStyleId = "{48338D13-A9CA-4A3B-89B7-9073ACB4FE37}"
$newStyleId = "{631E9890-4A1E-43CC-99DB-8398D616D4D7}"

$p = $ppt.ActivePresentation

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $s = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
